$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- New row of data on "strategy_id-1" (sheet2) ---
# Apply the black-font style to the whole new row range first (creates the
# new font/cellXf: font color rgb FF000000, same Calibri/11/family2/minor).
$ws2.Range("A2:AS2").Font.Color = 0

# Text columns (shared strings "General" / "frac_gnrl_eating_red_meat")
$ws2.Range("A2").Value = "General"
$ws2.Range("B2").Value = "frac_gnrl_eating_red_meat"

# C2:G2 stay blank (styled empty cells)

# Numeric trajectory values across H2:AS2
$ws2.Range("H2").Value = 1
$ws2.Range("I2").Value = 1
$ws2.Range("J2").Value = 1
$ws2.Range("K2").Value = 0.99
$ws2.Range("L2").Value = 0.96499999999999997
$ws2.Range("M2").Value = 0.94
$ws2.Range("N2").Value = 0.91500000000000004
$ws2.Range("O2").Value = 0.89
$ws2.Range("P2").Value = 0.86499999999999999
$ws2.Range("Q2").Value = 0.84
$ws2.Range("R2").Value = 0.81499999999999995
$ws2.Range("S2").Value = 0.79
$ws2.Range("T2").Value = 0.76500000000000001
$ws2.Range("U2").Value = 0.74
$ws2.Range("V2").Value = 0.71499999999999997
$ws2.Range("W2").Value = 0.69
$ws2.Range("X2").Value = 0.66500000000000004
$ws2.Range("Y2").Value = 0.64
$ws2.Range("Z2").Value = 0.61499999999999999
$ws2.Range("AA2").Value = 0.59
$ws2.Range("AB2").Value = 0.56499999999999995
$ws2.Range("AC2").Value = 0.54
$ws2.Range("AD2").Value = 0.51500000000000001
$ws2.Range("AE2").Value = 0.49
$ws2.Range("AF2").Value = 0.46500000000000002
$ws2.Range("AG2").Value = 0.44
$ws2.Range("AH2").Value = 0.41499999999999898
$ws2.Range("AI2").Value = 0.38999999999999901
$ws2.Range("AJ2").Value = 0.36499999999999899
$ws2.Range("AK2").Value = 0.33999999999999903
$ws2.Range("AL2").Value = 0.314999999999999
$ws2.Range("AM2").Value = 0.28999999999999898
$ws2.Range("AN2").Value = 0.26499999999999901
$ws2.Range("AO2").Value = 0.23999999999999899
$ws2.Range("AP2").Value = 0.214999999999999
$ws2.Range("AQ2").Value = 0.189999999999999
$ws2.Range("AR2").Value = 0.16499999999999901
$ws2.Range("AS2").Value = 0.13999999999999899

# --- Selection / active sheet changes ---
# Move sheet1's selection, then switch to and select on sheet2 so that
# sheet2 ("strategy_id-1") ends up as the active/visible tab.
[void]$ws1.Range("I10").Select()
[void]$ws2.Activate()
[void]$ws2.Range("F8").Select()
